$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold/border/center) from N1 onto the two new header cells
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# Header row (row 1): sequential index values 14, 15
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Data rows 2-67: new predicted values for columns O (t+14) and P (t+15)
$ws.Range("O2").Value = -0.2118667684557863
$ws.Range("P2").Value = -0.2109882148331458
$ws.Range("O3").Value = 0.2793020923570642
$ws.Range("P3").Value = 0.2790666899601203
$ws.Range("O4").Value = 0.2472161537314162
$ws.Range("P4").Value = 0.2493368564385359
$ws.Range("O5").Value = -0.03071681582589129
$ws.Range("P5").Value = -0.03223859973476871
$ws.Range("O6").Value = 0.2284998933090776
$ws.Range("P6").Value = 0.2276702295113165
$ws.Range("O7").Value = -0.3700144451782662
$ws.Range("P7").Value = -0.3688008954764657
$ws.Range("O8").Value = -0.1405647412431064
$ws.Range("P8").Value = -0.1362849091654254
$ws.Range("O9").Value = -0.282991153910005
$ws.Range("P9").Value = -0.2799835056006633
$ws.Range("O10").Value = 0.4291866612634972
$ws.Range("P10").Value = 0.4279199560978243
$ws.Range("O11").Value = -0.176344167883663
$ws.Range("P11").Value = -0.1760802441315965
$ws.Range("O12").Value = -0.009299640980374127
$ws.Range("P12").Value = -0.01446965669159919
$ws.Range("O13").Value = -0.02006580735222931
$ws.Range("P13").Value = -0.02232463524029228
$ws.Range("O14").Value = 0.1993789965724475
$ws.Range("P14").Value = 0.1907182614454261
$ws.Range("O15").Value = 0.102471931569647
$ws.Range("P15").Value = 0.0885226072407119
$ws.Range("O16").Value = 0.5498174582048032
$ws.Range("P16").Value = 0.5351402051065315
$ws.Range("O17").Value = 0.6297276653086581
$ws.Range("P17").Value = 0.613090013645827
$ws.Range("O18").Value = -0.07037916100497743
$ws.Range("P18").Value = -0.07859296945098027
$ws.Range("O19").Value = 0.4091996973908547
$ws.Range("P19").Value = 0.4004075141166942
$ws.Range("O20").Value = 0.4514442792950968
$ws.Range("P20").Value = 0.434540240673098
$ws.Range("O21").Value = 0.6674952515933892
$ws.Range("P21").Value = 0.6540509110686575
$ws.Range("O22").Value = 0.4369651510460372
$ws.Range("P22").Value = 0.4226486676943552
$ws.Range("O23").Value = -0.01514212134937037
$ws.Range("P23").Value = -0.02767112068261991
$ws.Range("O24").Value = 2.117460703324006
$ws.Range("P24").Value = 1.925728045253889
$ws.Range("O25").Value = 0.2945850229093873
$ws.Range("P25").Value = 0.2923396602738536
$ws.Range("O26").Value = 0.1586908415937264
$ws.Range("P26").Value = 0.1498131679641713
$ws.Range("O27").Value = 0.05307925167711211
$ws.Range("P27").Value = 0.04440415632521851
$ws.Range("O28").Value = 0.7962911922199494
$ws.Range("P28").Value = 0.7888547967758671
$ws.Range("O29").Value = 1.931507616094253
$ws.Range("P29").Value = 1.80345542120025
$ws.Range("O30").Value = 0.6337107333660277
$ws.Range("P30").Value = 0.6271640581275549
$ws.Range("O31").Value = -0.4976744190406188
$ws.Range("P31").Value = -0.4978915341616721
$ws.Range("O32").Value = 0.5303857589845878
$ws.Range("P32").Value = 0.5250016368730054
$ws.Range("O33").Value = 0.7377147283516161
$ws.Range("P33").Value = 0.7354447487251349
$ws.Range("O34").Value = -0.8532405706176763
$ws.Range("P34").Value = -0.8565720975500544
$ws.Range("O35").Value = 0.7759598044374806
$ws.Range("P35").Value = 0.777928761326932
$ws.Range("O36").Value = 0.7155991147306636
$ws.Range("P36").Value = 0.7200929221519748
$ws.Range("O37").Value = 0.6785997705589175
$ws.Range("P37").Value = 0.6830751538652481
$ws.Range("O38").Value = 0.6256015949527934
$ws.Range("P38").Value = 0.6259816841563289
$ws.Range("O39").Value = 0.5784468498126828
$ws.Range("P39").Value = 0.5808803528519939
$ws.Range("O40").Value = 0.7407471754940677
$ws.Range("P40").Value = 0.7427969830964934
$ws.Range("O41").Value = 0.5507243670676554
$ws.Range("P41").Value = 0.5536811331453142
$ws.Range("O42").Value = 0.5734184539660087
$ws.Range("P42").Value = 0.576323849082486
$ws.Range("O43").Value = 0.6516704305503593
$ws.Range("P43").Value = 0.6532973507272232
$ws.Range("O44").Value = 0.6664097011043051
$ws.Range("P44").Value = 0.6700023818368632
$ws.Range("O45").Value = 0.6260287946305059
$ws.Range("P45").Value = 0.6334119014180395
$ws.Range("O46").Value = -1.286871801810674
$ws.Range("P46").Value = -1.288973829779031
$ws.Range("O47").Value = -1.004966964109841
$ws.Range("P47").Value = -1.00650088910649
$ws.Range("O48").Value = -0.8787254424327713
$ws.Range("P48").Value = -0.8781591945888401
$ws.Range("O49").Value = -0.6409742573158563
$ws.Range("P49").Value = -0.6407390600474037
$ws.Range("O50").Value = -0.06068047698041935
$ws.Range("P50").Value = -0.06202798607491479
$ws.Range("O51").Value = -0.8660064416435256
$ws.Range("P51").Value = -0.8648277138693918
$ws.Range("O52").Value = -0.8660064416435256
$ws.Range("P52").Value = -0.8648277138693918
$ws.Range("O53").Value = -1.140135305052401
$ws.Range("P53").Value = -1.140506826939457
$ws.Range("O54").Value = -0.1745752278779049
$ws.Range("P54").Value = -0.174662191470809
$ws.Range("O55").Value = -1.023946891947809
$ws.Range("P55").Value = -1.025819070949864
$ws.Range("O56").Value = -0.8933621222011922
$ws.Range("P56").Value = -0.8963360751404963
$ws.Range("O57").Value = -0.9031413180167853
$ws.Range("P57").Value = -0.9087482464168485
$ws.Range("O58").Value = -1.048102919737715
$ws.Range("P58").Value = -1.0528720696564
$ws.Range("O59").Value = -0.7790363348415635
$ws.Range("P59").Value = -0.7796261854117366
$ws.Range("O60").Value = -0.4031864498251687
$ws.Range("P60").Value = -0.4049132571904611
$ws.Range("O61").Value = 0.3879701709778287
$ws.Range("P61").Value = 0.3878444624785629
$ws.Range("O62").Value = -1.15350588436076
$ws.Range("P62").Value = -1.159924479454534
$ws.Range("O63").Value = -0.5605882991002087
$ws.Range("P63").Value = -0.5568630763340242
$ws.Range("O64").Value = -0.8361923911266753
$ws.Range("P64").Value = -0.8363934034377949
$ws.Range("O65").Value = -0.02164126409552751
$ws.Range("P65").Value = -0.02279234689934571
$ws.Range("O66").Value = -0.7318678281875393
$ws.Range("P66").Value = -0.7384082241605199
$ws.Range("O67").Value = -0.7011768643852268
$ws.Range("P67").Value = -0.7101789725791851
